$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Rows.Item(2).Delete()
[void]$ws.Range("A2:K2").Select()
